# "My favorites testcase for osprey"
# Adds a new "MyFavorites" worksheet (after "GiftRegistry") mirroring the
# structure of the existing account-style sheets, plus small selection
# tweaks on two pre-existing sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) CreateAccount sheet: selection moves from A1:J6 to F17 (no tab switch
#    is left active on this sheet afterwards, so select another sheet
#    after this one).
# ---------------------------------------------------------------------
$wsCreateAccount = $wb.Worksheets.Item("CreateAccount")
$wsCreateAccount.Range("F17").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) GiftRegistry sheet: selection moves from J4 to A6, and it stops
#    being the tab shown when the workbook opens (topLeftCell="B1" is
#    also cleared along with the tab switch).
# ---------------------------------------------------------------------
$wsGiftRegistry = $wb.Worksheets.Item("GiftRegistry")
$wsGiftRegistry.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) New "MyFavorites" worksheet, inserted after "GiftRegistry".
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "MyFavorites"

# Header row (yellow-highlighted, matches the other account sheets).
$ws.Range("A1:O1").Interior.Color = 65535
$ws.Range("A1").Value = "DataSet"
$ws.Range("B1").Value = "UserName"
$ws.Range("C1").Value = "Prod UserName"
$ws.Range("D1").Value = "Password"
$ws.Range("E1").Value = "Confirm Password"
$ws.Range("F1").Value = "FirstName"
$ws.Range("G1").Value = "LastName"
$ws.Range("H1").Value = "Street"
$ws.Range("I1").Value = "City"
$ws.Range("J1").Value = "Region"
$ws.Range("K1").Value = "postcode"
$ws.Range("L1").Value = "phone"
$ws.Range("M1").Value = "Products"
$ws.Range("N1").Value = "Quantity"
$ws.Range("O1").Value = "methods"

# Row 2 - account + address details used for this testcase.
$ws.Range("A2").Value = "Account"

$ws.Range("B2").Value = "testersemail.278@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:testersemail.278@gmail.com") | Out-Null
$ws.Range("B2").Style = "Hyperlink"

$ws.Range("C2").Style = "Hyperlink"

$ws.Range("D2").Value = "Testers@278"
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:Testers@278") | Out-Null
$ws.Range("D2").Style = "Hyperlink"

$ws.Range("E2").Value = "Testers@278"
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Testers@278") | Out-Null
$ws.Range("E2").Style = "Hyperlink"

$ws.Range("F2").Value = "QA"
$ws.Range("G2").Value = "TEST"
$ws.Range("L2").Value = "'9898989898"
$ws.Range("I2").Value = "Little Rock"
$ws.Range("J2").Value = "Arkansas"
$ws.Range("K2").Value = "'72211"
$ws.Range("H2").Value = "6 Walnut Valley Dr"

# Row 3 - a favorited product.
$ws.Range("A3").Value = "Product"
$ws.Range("M3").Value = "POCO® CARRYING CASE"
$ws.Range("N3").Value = "'1"

# Row 4 - shipping method.
$ws.Range("A4").Value = "GroundShipping method"
$ws.Range("O4").Value = "Fixed - Flat Rate"

# Final selection / active tab: MyFavorites becomes the active sheet.
$ws.Range("E7").Select() | Out-Null
